# "fixed daily procents calculation"
#
# On the "Лист2" (sheet2) DepositDailyLine property list, insert a new
# field "decimal DepoRate" above the old "decimal DayProfit" row and
# rename the daily-profit fields to the more accurate "Procents" naming:
#   decimal Balance
#   decimal DepoRate          <- new
#   decimal DayProcents       <- renamed from "decimal DayProfit"
#   decimal NotPaidProcents   <- renamed from "decimal NotPaidProfit"
#   decimal DayDevaluation

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист2")

# Push "decimal DayProfit" / "decimal NotPaidProfit" / "decimal DayDevaluation"
# (and everything below them) down by one row to make room for the new field.
$ws.Rows(26).Insert()

$ws.Range("D26").Value = "decimal DepoRate"
$ws.Range("D27").Value = "decimal DayProcents"
$ws.Range("D28").Value = "decimal NotPaidProcents"

# Match the author's final cursor position / view after the edit.
$ws.Activate()
$ws.Range("D29").Select()
